$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the credential values in row 3 while leaving row 1/2 untouched.
$ws.Range("A3").Value = "s6670405"
$ws.Range("B3").Value = "Quality77"
